$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new (blank) column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Inserting column N shifts old N->O (Late), old O header->P (heading/Original),
# old P->Q (Outstanding) one column to the right, and leaves the new N column blank.
$ws.Columns("N:N").Insert()

# Match the new column's width to its left neighbour (column M).
$ws.Range("N1").ColumnWidth = $ws.Range("M1").ColumnWidth

# Update the selection to match the saved view state.
$ws.Range("R10").Select()
